# New crime data collected - update CompStat weekly report
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
# "Volume 31   Number  45" -> "...  46"
$ws.Range("A8").Characters(21, 2).Text = "46"
# "Report Covering the Week  11/4/2024  Through  11/10/2024"
#   -> "...  11/11/2024  Through  11/17/2024"
$ws.Range("C9").Characters(27, 9).Text = "11/11/2024"
$ws.Range("C9").Characters(48, 10).Text = "11/17/2024"

# --- Column width tweak (col I / 9) ---
# Target stored width is 7.433768 (same as column H). The engine's
# ColumnWidth setter re-adds its own ~5/7 padding internally, so feed it
# the value net of that padding to land on the closest achievable width.
$ws.Columns.Item(9).ColumnWidth = (7.433768 - 5/7)

# --- Row 14 (Murder) ---
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("C14").Value = 1
$ws.Range("F14").Value = 2
$ws.Range("I14").Value = 6
$ws.Range("K14").Value = 100
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -71.428571428571

# --- Row 15 (Rape) ---
$ws.Range("G15").Value = 2
$ws.Range("N15").Value = -86.666666666666

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 71.428571428571
$ws.Range("I16").Value = 90
$ws.Range("J16").Value = 86
$ws.Range("K16").Value = 4.651162790697
$ws.Range("L16").Value = -8.163265306122
$ws.Range("M16").Value = -57.142857142857
$ws.Range("N16").Value = -91.295938104448

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -11.764705882352
$ws.Range("I17").Value = 205
$ws.Range("J17").Value = 157
$ws.Range("K17").Value = 30.573248407643
$ws.Range("L17").Value = 36.666666666666
$ws.Range("M17").Value = 75.213675213675
$ws.Range("N17").Value = -53.514739229024

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 4
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 9
$ws.Range("H18").Value = 28.571428571428
$ws.Range("I18").Value = 123
$ws.Range("J18").Value = 87
$ws.Range("K18").Value = 41.379310344827
$ws.Range("L18").Value = 12.844036697247
$ws.Range("M18").Value = -22.641509433962
$ws.Range("N18").Value = -87.253886010362

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 19.230769230769
$ws.Range("I19").Value = 257
$ws.Range("J19").Value = 285
$ws.Range("K19").Value = -9.824561403508
$ws.Range("L19").Value = -16.558441558441
$ws.Range("M19").Value = -37.621359223301
$ws.Range("N19").Value = -47.764227642276

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 11
$ws.Range("H20").Value = 57.142857142857
$ws.Range("I20").Value = 116
$ws.Range("J20").Value = 88
$ws.Range("K20").Value = 31.818181818181
$ws.Range("L20").Value = 39.759036144578
$ws.Range("M20").Value = 31.818181818181
$ws.Range("N20").Value = -82.262996941896

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 37.5
$ws.Range("F21").Value = 80
$ws.Range("G21").Value = 66
$ws.Range("H21").Value = 21.212121212121
$ws.Range("I21").Value = 801
$ws.Range("J21").Value = 716
$ws.Range("K21").Value = 11.871508379888
$ws.Range("L21").Value = 5.533596837944
$ws.Range("M21").Value = -19.9
$ws.Range("N21").Value = -77.976354138025

# --- Row 22 (Transit) ---
$ws.Range("D22").Copy($ws.Range("C22"))
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -72.549019607843

# --- Row 23 (Housing) ---
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -25
$ws.Range("I23").Value = 92
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = -8
$ws.Range("L23").Value = 19.480519480519
$ws.Range("M23").Value = 2.222222222222

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = -26.666666666666
$ws.Range("F24").Value = 41
$ws.Range("G24").Value = 64
$ws.Range("H24").Value = -35.9375
$ws.Range("I24").Value = 518
$ws.Range("J24").Value = 681
$ws.Range("K24").Value = -23.935389133627
$ws.Range("L24").Value = -23.598820058997
$ws.Range("M24").Value = -49.463414634146

# --- Row 25 (Retail Theft) ---
$ws.Range("C25").Value = 1
$ws.Range("E25").Value = -50
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = -42.857142857142
$ws.Range("J25").Value = 159
$ws.Range("K25").Value = -20.125786163522
$ws.Range("L25").Value = -30.978260869565

# --- Row 26 (Misd. Assault) ---
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 200
$ws.Range("F26").Value = 33
$ws.Range("G26").Value = 23
$ws.Range("H26").Value = 43.478260869565
$ws.Range("I26").Value = 342
$ws.Range("J26").Value = 297
$ws.Range("K26").Value = 15.151515151515
$ws.Range("L26").Value = 1.483679525222
$ws.Range("M26").Value = 32.046332046332

# --- Row 27 (UCR Rape*) ---
$ws.Range("C27").Copy($ws.Range("D27"))
$ws.Range("M27").Copy($ws.Range("E27"))
$ws.Range("G27").Value = 3

# --- Row 28 (Other Sex Crimes) ---
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 29
$ws.Range("K28").Value = -27.5
$ws.Range("L28").Value = 11.538461538461

# --- Row 29 (Shooting Vic.) ---
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 10
$ws.Range("K29").Value = 30
$ws.Range("N29").Value = -88.596491228070

# --- Row 30 (Shooting Inc.) ---
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 9
$ws.Range("K30").Value = 0
$ws.Range("N30").Value = -91

# --- Row 31 (Hate Crimes) ---
$ws.Range("C31").Copy($ws.Range("D31"))
$ws.Range("M31").Copy($ws.Range("E31"))
$ws.Range("F31").Value = 4
$ws.Range("H31").Value = 100
$ws.Range("I31").Value = 5
$ws.Range("K31").Value = -16.666666666666
$ws.Range("L31").Value = -16.666666666666
